$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# H1: 15:00 - 17:00 -> 15:00 - 23:00
$ws.Range("H1").Value = "15:00 - 23:00"

# Row 4 (Wednesday): E4, F4, G4 -> TEST, and add H4 = TEST
$ws.Range("E4").Value = "TEST"
$ws.Range("F4").Value = "TEST"
$ws.Range("G4").Value = "TEST"
$ws.Range("H4").Value = "TEST"
